$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.6357677902621723
$ws1.Range("C2").Value = 0.5787187839305103
$ws1.Range("D2").Value = 0.99812734082397
$ws1.Range("E2").Value = 0.7326460481099656
$ws1.Range("F2").Value = 0.8717697088649002
$ws1.Range("G2").Value = 0.9710601919977577
$ws1.Range("H2").Value = 0.8023502924714893
$ws1.Range("I2").Value = 533
$ws1.Range("J2").Value = 388
$ws1.Range("K2").Value = 146
$ws1.Range("L2").Value = 1

# --- Sheet: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")

$ws2.Range("B2").Value = 0.9931972789115646
$ws2.Range("C2").Value = 0.2734082397003745
$ws2.Range("D2").Value = 0.4287812041116006

$ws2.Range("B3").Value = 0.5787187839305103
$ws2.Range("C3").Value = 0.99812734082397
$ws2.Range("D3").Value = 0.7326460481099656

$ws2.Range("B4").Value = 0.6357677902621723
$ws2.Range("C4").Value = 0.6357677902621723
$ws2.Range("D4").Value = 0.6357677902621723
$ws2.Range("E4").Value = 0.6357677902621723

$ws2.Range("B5").Value = 0.7859580314210375
$ws2.Range("C5").Value = 0.6357677902621723
$ws2.Range("D5").Value = 0.580713626110783

$ws2.Range("B6").Value = 0.7859580314210375
$ws2.Range("C6").Value = 0.6357677902621723
$ws2.Range("D6").Value = 0.5807136261107831

# --- Sheet: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 146
$ws3.Range("C2").Value = 388
$ws3.Range("B3").Value = 1
$ws3.Range("C3").Value = 533
